# Procesos_AF.xlsx — "Add files via upload"
#
# The sheet that used to hold the C_11 series now carries the C_16.1
# series, so rename the (single) worksheet and let the workbook-level
# _FilterDatabase defined name follow it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "C_16.1"

# Excel quotes sheet names that contain a "." when it rewrites a
# reference, e.g. 'C_16.1'!$B$5:$B$5 — make sure the hidden
# _xlnm._FilterDatabase name (built from the old AutoFilter range on
# B5) is refreshed to match explicitly, in case the rename alone
# doesn't normalize the quoting.
$filterName = $wb.Names.Item("_xlnm._FilterDatabase")
$filterName.RefersTo = "='C_16.1'!`$B`$5:`$B`$5"
